$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 53818
$ws.Range("B2").Value = "Dr. Thales Martins"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45088
$ws.Range("G2").Value = 7830.25

# Row 3
$ws.Range("A3").Value = 81255
$ws.Range("B3").Value = "Lorenzo Barros"
$ws.Range("C3").Value = "Jurídico"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4470.51

# Row 4
$ws.Range("A4").Value = 44294
$ws.Range("B4").Value = "Eduardo Ramos"
$ws.Range("D4").Value = "Doença"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 10596.73

# Row 5
$ws.Range("A5").Value = 67602
$ws.Range("B5").Value = "Heitor Duarte"
$ws.Range("C5").Value = "P&D"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 45093
$ws.Range("G5").Value = 5301.81

# Row 6
$ws.Range("A6").Value = 95467
$ws.Range("B6").Value = "Sr. Vitor Gabriel Gomes"
$ws.Range("C6").Value = "Recursos Humanos"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45092
$ws.Range("G6").Value = 8625.219999999999

# Row 7
$ws.Range("A7").Value = 33511
$ws.Range("B7").Value = "Letícia Freitas"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 45089
$ws.Range("G7").Value = 3130.36

# Row 8
$ws.Range("A8").Value = 81521
$ws.Range("B8").Value = "Davi Luiz Duarte"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45081
$ws.Range("G8").Value = 10348.93

# Row 9
$ws.Range("A9").Value = 99725
$ws.Range("B9").Value = "Kevin da Cunha"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45104
$ws.Range("G9").Value = 12409.62

# Row 10
$ws.Range("A10").Value = 4925
$ws.Range("B10").Value = "Sra. Raquel Fernandes"
$ws.Range("C10").Value = "Vendas"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("F10").Value = 45083
$ws.Range("G10").Value = 12063.29

# Row 11
$ws.Range("A11").Value = 84659
$ws.Range("B11").Value = "Luigi Rezende"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45103
$ws.Range("G11").Value = 8516.43
